$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: add new column I, shift G1/H1 values ---
# Give the new I1 header cell the same style (bold, centered, bordered) as the other header cells
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 10

# --- Data rows 2-10: updated p-values across columns B-I ---
# Row 2
$ws.Range("B2").Value = 0.000000003313112006964047
$ws.Range("C2").Value = 0.000000040082434704658
$ws.Range("D2").Value = 0.000000000001692868067948439
$ws.Range("E2").Value = 0.00000000000001176836406102666
$ws.Range("F2").Value = 0.0000000000002757793993168889
$ws.Range("G2").Value = 0.000000000001337374655463464
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.00000000000002819966482547898

# Row 3
$ws.Range("B3").Value = 0.0000285649688822609
$ws.Range("C3").Value = 0.000000006161926524583805
$ws.Range("D3").Value = 0.00000000001567812546454661
$ws.Range("E3").Value = 0.00000000000000333066907387547
$ws.Range("F3").Value = 0.0000000000005255795798575491
$ws.Range("G3").Value = 0.000000000003168132423070347
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.00000000000004019007349143067

# Row 4
$ws.Range("B4").Value = 0.000000001384601544529573
$ws.Range("C4").Value = 0.00000004280528442812681
$ws.Range("D4").Value = 0.000000000001389111048410996
$ws.Range("E4").Value = 0.0000000000000113242748511766
$ws.Range("F4").Value = 0.000000000000269784194983913
$ws.Range("G4").Value = 0.000000000001165956220461339
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.00000000000002775557561562891

# Row 5
$ws.Range("B5").Value = 0.008512010379494672
$ws.Range("C5").Value = 0.000000003285452576662351
$ws.Range("D5").Value = 0.000000000005679012815562601
$ws.Range("E5").Value = 0.0000000000004514166818125886
$ws.Range("F5").Value = 0.0000000000006008527009271347
$ws.Range("G5").Value = 0.000000000005044409334686861
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.0000000000000184297022087776

# Row 6
$ws.Range("B6").Value = 0.000000000000007327471962526033
$ws.Range("C6").Value = 0.0000001059631213884416
$ws.Range("D6").Value = 0.00000000001344813149728452
$ws.Range("E6").Value = 0.0000000000000004440892098500626
$ws.Range("F6").Value = 0.00000000001069877519910278
$ws.Range("G6").Value = 0.000000000008893108471852429
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0.00000000000007838174553853605

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.00000002271116272822837
$ws.Range("D7").Value = 0.000000000003030686812621752
$ws.Range("E7").Value = 0.00000000000002020605904817785
$ws.Range("F7").Value = 0.000000000001209254918421721
$ws.Range("G7").Value = 0.000000000004530154029680489
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.000000000000008215650382226158

# Row 8
$ws.Range("B8").Value = 0.00000000000001998401444325282
$ws.Range("C8").Value = 0.00000002973109003434615
$ws.Range("D8").Value = 0.0000000000134878774815661
$ws.Range("E8").Value = 0.00000000000003774758283725532
$ws.Range("F8").Value = 0.000000000000595967719618784
$ws.Range("G8").Value = 0.000000000003734568210234102
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.00000000000004862776847858186

# Row 9
$ws.Range("B9").Value = 0.000002547181449541469
$ws.Range("C9").Value = 0.00000005560486515676644
$ws.Range("D9").Value = 0.000000000005720979245893432
$ws.Range("E9").Value = 0.00000000000002375877272697835
$ws.Range("F9").Value = 0.0000000000004505285033928885
$ws.Range("G9").Value = 0.0000000000023046009545169
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0.0000000000000297539770599542

# Row 10
$ws.Range("B10").Value = 0.6067316732500676
$ws.Range("C10").Value = 0.4022439599054577
$ws.Range("D10").Value = 0.9856428175419636
$ws.Range("E10").Value = 0.9066100980529868
$ws.Range("F10").Value = 0.06755822882963347
$ws.Range("G10").Value = 0.00000000002050182246193799
$ws.Range("H10").Value = 0.0000000002186013592364588
$ws.Range("I10").Value = 0.000000000001362021606610142

